$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "all" (sheet1.xml): append a new daily row (was last row = footnote
# at row 40), pushing the footnote down to row 41.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

# Move the footnote row (currently B40) down to B41, carrying its style.
$wsAll.Range("B40").Copy($wsAll.Range("B41"))

# Use row 39 as a style template for the new data row 40.
$wsAll.Range("A39:H39").Copy($wsAll.Range("A40:H40"))

# Now overwrite the new row 40 with the new day's figures.
$wsAll.Range("A40").Value = 43968
$wsAll.Range("B40").Value = 283
$wsAll.Range("C40").Value = 280
$wsAll.Range("D40").Value = 51
$wsAll.Range("E40").Value = 44
$wsAll.Range("F40").Value = 7
$wsAll.Range("G40").Value = 11
$wsAll.Range("H40").Value = 218

# ---------------------------------------------------------------------------
# Sheet "kobe" (sheet2.xml): correct the last data row (94), append a new
# daily row (95), pushing the footnote down to row 96.
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

# Correct already-published figures for row 94 (revision of 5/16 data).
$wsKobe.Range("D94").Value = 1
$wsKobe.Range("E94").Value = 283
$wsKobe.Range("F94").Value = 47
$wsKobe.Range("G94").Value = 41
$wsKobe.Range("J94").Value = 208

# Move the footnote row (currently B95) down to B96, carrying its style.
$wsKobe.Range("B95").Copy($wsKobe.Range("B96"))

# Use row 94 as a style template for the new data row 95.
$wsKobe.Range("A94:J94").Copy($wsKobe.Range("A95:J95"))

# Now overwrite the new row 95 with the new day's figures.
$wsKobe.Range("A95").Value = 43968
$wsKobe.Range("B95").Value = 25
$wsKobe.Range("C95").Value = 2865
$wsKobe.Range("D95").Value = 0
$wsKobe.Range("E95").Value = 283
$wsKobe.Range("F95").Value = 46
$wsKobe.Range("G95").Value = 40
$wsKobe.Range("H95").Value = 6
$wsKobe.Range("I95").Value = 11
$wsKobe.Range("J95").Value = 209

# ---------------------------------------------------------------------------
# Update the "kobe" sheet's own remembered selection (bottomRight pane now
# points at the new last row), without leaving it as the active tab.
# ---------------------------------------------------------------------------
$wsKobe.Activate()
$wsKobe.Range("A96").Select()

# ---------------------------------------------------------------------------
# Tab selection: move the active tab from "other" to "all" (leaves "other"'s
# own selection/pane state untouched).
# ---------------------------------------------------------------------------
$wsAll.Activate()
$wsAll.Range("B41").Select()
